$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the new JSON-formatted questions string (matches json.dumps(questions, indent=4))
# Note: json.dumps with ensure_ascii=True (the default) escapes the right single
# quotation mark (U+2019) as the *literal* six characters \u2019 rather than
# emitting the actual Unicode character, so we reproduce that literal text here.
$rightSingleQuote = "\u2019"

$newText = "questions = [" + "`n" +
"    {" + "`n" +
"        `"title`": `"Your company wants to determine the total cost (y) function as a function of the number of units produced (x). In this regard, the following estimates were provided by the company" + $rightSingleQuote + "s accountant:Fixed expenditure for the year of `$100,000Raw material cost for each unit produced of `$4Labor cost per unit of `$1Packing and shipping cost per unit of `$2  Which of the following best describes the total cost (y) function of your company?`"," + "`n" +
"        `"ques_type`": 2," + "`n" +
"        `"options`": [" + "`n" +
"            `" y = 4x + x + 2x + 100,000`"," + "`n" +
"            `"y + 100,000 = 4x + x + 2x`"," + "`n" +
"            `"y = 4x + 100,000 (x+2x)`"," + "`n" +
"            `"y + 100,000 (x+2x) = 4x`"" + "`n" +
"        ]," + "`n" +
"        `"score`": `"y = 4x + x + 2x + 100,000`"" + "`n" +
"    }," + "`n" +
"    {" + "`n" +
"        `"title`": `"Your company has limited materials, as shown in the table below, to produce Product-A and Product-B. You need to determine the quantity of these products (the only ones your company manufactures) to maximize profits. You created a graphical representation and determined the optimal solution exists at point B, where lines for these two constraints intersect with each other. How many units of each product should be produced to maximize profit?`"," + "`n" +
"        `"ques_type`": 2," + "`n" +
"        `"options`": [" + "`n" +
"            `"40 units of A and 90 units of B`"," + "`n" +
"            `"60 units of A and 100 units of B`"," + "`n" +
"            `"90 units of A and 40 units of B`"," + "`n" +
"            `"80 units of A and 100 units of B`"" + "`n" +
"        ]," + "`n" +
"        `"score`": `"90 units of A and 40 units of B`"" + "`n" +
"    }," + "`n" +
"    {" + "`n" +
"        `"title`": `"Your company" + $rightSingleQuote + "s bank offered placement of funds at a 10% annual rate of return, compounded monthly.  What is the effective annual rate of return of these funds?`"," + "`n" +
"        `"ques_type`": 2," + "`n" +
"        `"options`": [" + "`n" +
"            `"10.00%`"," + "`n" +
"            `"10.47%`"," + "`n" +
"            `"11.55%`"," + "`n" +
"            `"12.00%`"" + "`n" +
"        ]," + "`n" +
"        `"score`": `"10.47%`"" + "`n" +
"    }," + "`n" +
"    {" + "`n" +
"        `"title`": `"Your company is considering investing `$50,000 in a five-year project. This project expects annual cash inflows of `$20,000 for the first four years and cash inflows of `$30,000 for the fifth year. The initial investment is expected to have a residual value of `$10,000 at the end of the fifth year, and the required rate of return is 10%. What is the present value of cash inflows associated with the project?`"," + "`n" +
"        `"ques_type`": 2," + "`n" +
"        `"options`": [" + "`n" +
"            `"`$78,234`"," + "`n" +
"            `"`$82,025`"," + "`n" +
"            `"`$88,234`"," + "`n" +
"            `"`$92,025`"" + "`n" +
"        ]," + "`n" +
"        `"score`": `"`$88,234`"" + "`n" +
"    }" + "`n" +
"]"

# Clear the old contents (A1 and A2) and remove any styling/borders applied to A1
$ws.Cells.Item(1, 1).Value = $null
$ws.Cells.Item(2, 1).Value = $null
$ws.Cells.Item(1, 1).ClearFormats()
$ws.Cells.Item(2, 1).ClearFormats()

# Write the new consolidated JSON text into A1 only
$ws.Range("A1").Value = $newText

# The multi-line text would otherwise trigger an automatic custom row
# height; AutoFit() recomputes it back to the sheet's default so no stray
# ht/customHeight attributes are left behind on the row.
$ws.Rows.Item(1).AutoFit()
